# Add two new data rows (7 and 8) to Sheet1, matching the existing
# wt% log table layout (columns A..FI). Column A holds a text label;
# columns B..FI hold the numeric measurements for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 6, A6:FI6) down
# into the two new rows first, so the new rows inherit the same cell
# styles (e.g. the bold/bordered/centered style used on column A) as the
# rest of the table.
$ws.Range("A6:FI6").Copy($ws.Range("A7:FI7"))
$ws.Range("A6:FI6").Copy($ws.Range("A8:FI8"))

# ---- Row 7 ----
$ws.Range("A7").Value = "a1"

$row7Values = @(
    "0","0.1777111134880425","11.48760016186961","21.9307360896316","14.93593020200784","15.38416822324578","1.061269866106982","34.95694954363999",
    "0.05573916454090103","0.009849304206473728","0","1.899029989640977e-05","1.091718408562078e-06","1.801831473116669e-09","9.444529836483061e-15","2.77281392709501e-15",
    "2.83235049226907e-17","1.754762363228676e-18","8.911887477111833e-19","4.58672093908158e-22","6.109640544949654e-20","5.590642726103323e-35","2.03522404727945e-22","7.697633459060879e-12",
    "5.332008738360454e-32","1.756369956804511e-35","8.199768336057133e-10","0","4.391754106592498e-10","1.082505281506059e-27","1.632656868566204e-11","6.098588748925914e-13",
    "1.87050626250231e-05","7.541094327229084e-06","1.85462599449264e-12","0","1.121938550287161e-14","6.187067989783687e-16","6.048974346833257e-19","0",
    "0","0","0.0004004358735331036","0.1461237947179045","0.7107868741520554","1.325639045890336","3.028839548874425","0.7136079534610805",
    "93.93601783077642","0.1026049637691531","0.0357316000546213","0","0.0001379311509387722","1.409710448199717e-05","2.816653288824494e-08","1.923080964826021e-13",
    "5.845141598008095e-14","9.727888470743817e-16","6.782624933746287e-17","3.375607136486949e-17","2.366284270772045e-20","3.194847474648148e-18","5.004001341463816e-33","9.983527535475029e-21",
    "4.106948938238706e-11","9.072840987727206e-31","1.294512260020593e-33","2.641532855840986e-09","1.688819698998617e-34","3.900773428518689e-10","7.72143718507519e-26","2.048642294488102e-11",
    "1.393803736517609e-12","6.252090440861975e-05","3.337200078873685e-05","8.450435185903895e-12","0","3.468617817573517e-14","3.523537585076893e-15","3.285437996252342e-18",
    "0","0","0","1.349129472782425e-07","3.393515664660708e-05","0.0001770907265838275","0.0003651882341929542","0.001991682722635607",
    "0.01380241839673461","99.91563199754771","0.0196520847249126","0.04609349092167735","0","0.001186107800365999","0.001064758816064875","1.109929724160228e-06",
    "7.773899147586279e-11","2.819811480322898e-11","3.311455340869854e-12","3.344621332939921e-13","2.210677089745808e-13","3.319041238037791e-16","1.252042638694171e-14","2.321536490148479e-24",
    "1.357373823897102e-16","2.121145550949315e-27","4.771536198332376e-24","9.609806272473323e-25","2.920034048323175e-28","3.763827459983333e-25","1.136082604560355e-33","6.419143216165825e-21",
    "1.159667875580017e-14","6.362070922984508e-16","2.697507342878574e-28","4.447290847643338e-28","2.1103712016158e-15","0","9.854447777561916e-18","6.414756307769566e-19",
    "6.322118581650009e-22","0","0","0","2.312970935849128e-10","1.397760684206232e-08","3.601762753504974e-08","4.122779256701957e-27",
    "1.606418342046371e-19","2.144440654546462e-12","0.004717505962711335","6.212491880027871e-07","0.0001445214588776984","0","0.01208685800985953","94.79104669315728",
    "0.01327951317626226","0.01120318715154396","0.009216398825507625","2.303875389944794","1.056108709514717","1.625116108506117","0.09672458639544869","0.02016962252428043",
    "2.759876602030989e-08","0.04465067781575134","5.009473588159614e-09","6.913081273855949e-07","2.579472480595032e-07","2.189796896624471e-07","1.132943243320096e-07","1.338106772421228e-07",
    "0.003614023446332803","2.918151146036188e-09","1.721572715939265e-10","0.006411953046805428","0.001632128178195646","3.649273194546234e-10","0","3.883839317906849e-12",
    "2.671952643580955e-13","2.635712709825498e-16","0","0"
)

for ($i = 0; $i -lt $row7Values.Length; $i++) {
    # $row7Values[0] goes into column B (index 2), etc.
    $ws.Cells.Item(7, $i + 2).Value = [double]$row7Values[$i]
}

# ---- Row 8 ----
$ws.Range("A8").Value = "b2"

$row8Values = @(
    "0","0.1777137668108373","11.48775716530647","21.93101807305863","14.93608719296298","15.38435511989394","1.061294706664068","34.95614080003494",
    "0.05573779274889248","0.00984901193810489","0","1.898617074919539e-05","1.091696917662313e-06","1.802624009389222e-09","9.448691440794087e-15","2.774116832709997e-15",
    "2.833351473880494e-17","1.755373460252352e-18","8.914876941114068e-19","4.588241698678974e-22","6.112245259024381e-20","5.587930103113365e-35","2.036140113282953e-22","7.7094156317721e-12",
    "5.328501128342067e-32","1.755815327665941e-35","8.208790279234165e-10","0","4.398483367124794e-10","1.082952653724485e-27","1.634554776093537e-11","6.107326608376492e-13",
    "1.873450558940313e-05","7.555117984748992e-06","1.856476872871014e-12","0","1.123192020679536e-14","6.193079989903404e-16","6.052095598413135e-19","0",
    "0","0","0.0004005163697291283","0.1461204517980234","0.7106670556005033","1.325309133388004","3.02788466623044","0.7134511659706533",
    "93.93758996637906","0.1025986885508421","0.03573037068423748","0","0.0001378390499630711","1.409162074651334e-05","2.816659991452854e-08","1.923175957250876e-13",
    "5.845888967367822e-14","9.727394910509877e-16","6.782542341547924e-17","3.375577062744559e-17","2.366095927906661e-20","3.194836339100556e-18","4.979223089559195e-33","9.983734302670209e-21",
    "4.110945593282527e-11","9.061152558495575e-31","1.293638357611616e-33","2.644048736577163e-09","1.687811453079516e-34","3.905212015745283e-10","7.722066037288865e-26","2.050802909106022e-11",
    "1.395322995015806e-12","6.26048057092253e-05","3.341827916940924e-05","8.455240210421312e-12","0","3.472095850945477e-14","3.525969217868057e-15","3.286983517536641e-18",
    "0","0","0","1.349239989249512e-07","3.393800728076814e-05","0.0001771057124425422","0.0003652205265354862","0.001991699571000888",
    "0.01380237671484023","99.91563246277923","0.01965208831656726","0.04609351600139634","0","0.001185703285552279","0.001064643787400637","1.110263860712835e-06",
    "7.776808026695186e-11","2.820849390757783e-11","3.312355431898541e-12","3.345486721096017e-13","2.211251985040121e-13","3.319849977235498e-16","1.252464719247464e-14","2.320964479291265e-24",
    "1.357853498456423e-16","2.121034234056981e-27","4.770541575299105e-24","9.608783226480859e-25","2.919062003473134e-28","3.763940016548926e-25","1.130920095345712e-33","6.421566074364222e-21",
    "1.161155474807687e-14","6.371829920328948e-16","2.69742945828522e-28","4.449891776910562e-28","2.112759416474628e-15","0","9.864938795407747e-18","6.422923759386954e-19",
    "6.326032822280702e-22","0","0","0","2.317156888839611e-10","1.400759849667456e-08","3.609360537182115e-08","4.122401796750682e-27",
    "1.606392848243687e-19","2.144621626918828e-12","0.004718713057963154","6.213394848796526e-07","0.0001445578636462119","0","0.01208459526624575","94.78927176192556",
    "0.01328216569071166","0.01120701702391937","0.009219552823044404","2.304664461664362","1.056470266127663","1.62567250018647","0.09675768600676452","0.02017653039355025",
    "2.759439804571768e-08","0.04466597123161932","5.009754507339577e-09","6.912318877191395e-07","2.579501408375656e-07","2.18978234267151e-07","1.13304445868688e-07","1.338173713141619e-07",
    "0.003615261294664294","2.918199161907355e-09","1.72168958018732e-10","0.006414149221096406","0.001632687202547547","3.648594467708288e-10","0","3.883508587408862e-12",
    "2.671747077503484e-13","2.633931113307559e-16","0","0"
)

for ($i = 0; $i -lt $row8Values.Length; $i++) {
    $ws.Cells.Item(8, $i + 2).Value = [double]$row8Values[$i]
}
